$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item('展览')
$ws1.Range('F4').Value = 5173
$ws1.Range('F5').Value = 5173
$ws1.Range('F6').Value = 122
$ws1.Range('F8').Value = 211
$ws1.Range('F10').Value = 213
$ws1.Range('F11').Value = 176
$ws1.Range('F12').Value = 8622
$ws1.Range('F13').Value = 8622
$ws1.Range('F14').Value = 275
$ws1.Range('F16').Value = 628
$ws1.Range('F17').Value = 4
$ws1.Range('F18').Value = 2572
$ws1.Range('F26').Value = 6501
$ws1.Range('F27').Value = 203
$ws1.Range('F29').Value = 144
$ws1.Range('F32').Value = 7027
$ws1.Range('F38').Value = 109
$ws1.Range('B40').NumberFormat = '@'
$ws1.Range('B40').Value = '2024-05-12'
$ws1.Range('C40').NumberFormat = '@'
$ws1.Range('C40').Value = '北京·次元仙界会'
$ws1.Range('D40').NumberFormat = '@'
$ws1.Range('D40').Value = '丽泽天地购物中心 丽泽天地购物中心'
$ws1.Range('E40').NumberFormat = '@'
$ws1.Range('E40').Value = '2024.05.12 10:00-05.13 02:00'
$ws1.Range('F40').Value = 0
$ws1.Range('G40').Value = 29
$ws1.Range('H40').NumberFormat = '@'
$ws1.Range('H40').Value = 'https://show.bilibili.com/platform/detail.html?id=83690'
$ws1.Range('I40').NumberFormat = '@'
$ws1.Range('I40').Value = '//i0.hdslb.com/bfs/openplatform/202404/nw4FcqlD1712029131170.png'
$ws1.Range('C41').NumberFormat = '@'
$ws1.Range('C41').Value = '北京·ICOS SP漫展04X五只猫动漫节'
$ws1.Range('D41').NumberFormat = '@'
$ws1.Range('D41').Value = '北京电影学院影视文化产业创新园平房园区 北京五只猫娱乐Mall'
$ws1.Range('E41').NumberFormat = '@'
$ws1.Range('E41').Value = '2024.05.18 09:00-05.19 17:00'
$ws1.Range('F41').Value = 39
$ws1.Range('G41').Value = 80
$ws1.Range('H41').NumberFormat = '@'
$ws1.Range('H41').Value = 'https://show.bilibili.com/platform/detail.html?id=83122'
$ws1.Range('I41').NumberFormat = '@'
$ws1.Range('I41').Value = '//i0.hdslb.com/bfs/openplatform/202403/3N8tBAKl1710831573887.jpeg'
$ws1.Range('C42').NumberFormat = '@'
$ws1.Range('C42').Value = '北京·YIYOU二次元大聚会'
$ws1.Range('D42').NumberFormat = '@'
$ws1.Range('D42').Value = '京开高速入口与京开高速交叉口西180米 北京双马文体创业园'
$ws1.Range('E42').NumberFormat = '@'
$ws1.Range('E42').Value = '2024.05.18 10:00-05.18 18:00'
$ws1.Range('F42').Value = 55
$ws1.Range('G42').Value = 55
$ws1.Range('H42').NumberFormat = '@'
$ws1.Range('H42').Value = 'https://show.bilibili.com/platform/detail.html?id=83129'
$ws1.Range('I42').NumberFormat = '@'
$ws1.Range('I42').Value = '//i1.hdslb.com/bfs/openplatform/202403/ZhTtVA3A1710812150528.png'
$ws1.Range('C43').NumberFormat = '@'
$ws1.Range('C43').Value = '北京·原神only3.0'
$ws1.Range('E43').NumberFormat = '@'
$ws1.Range('E43').Value = '2024.05.18 10:00-05.19 17:00'
$ws1.Range('F43').Value = 2542
$ws1.Range('G43').Value = 68
$ws1.Range('H43').NumberFormat = '@'
$ws1.Range('H43').Value = 'https://show.bilibili.com/platform/detail.html?id=81766'
$ws1.Range('I43').NumberFormat = '@'
$ws1.Range('I43').Value = '//i2.hdslb.com/bfs/openplatform/202402/Lfxwe5PO1707120983684.jpeg'
$ws1.Range('C44').NumberFormat = '@'
$ws1.Range('C44').Value = '北京·原神only3.0——32D小神奈签售会'
$ws1.Range('F44').Value = 36
$ws1.Range('H44').NumberFormat = '@'
$ws1.Range('H44').Value = 'https://show.bilibili.com/platform/detail.html?id=82147'
$ws1.Range('I44').NumberFormat = '@'
$ws1.Range('I44').Value = '//i1.hdslb.com/bfs/openplatform/202402/lQoExxJd1709100610683.jpeg'
$ws1.Range('C45').NumberFormat = '@'
$ws1.Range('C45').Value = '北京·原神only3.0——蛋黄mayo签售会'
$ws1.Range('D45').NumberFormat = '@'
$ws1.Range('D45').Value = '北花园路1号 超级蜂巢'
$ws1.Range('E45').NumberFormat = '@'
$ws1.Range('E45').Value = '2024.05.18 10:00-05.18 17:00'
$ws1.Range('F45').Value = 71
$ws1.Range('G45').Value = 1
$ws1.Range('H45').NumberFormat = '@'
$ws1.Range('H45').Value = 'https://show.bilibili.com/platform/detail.html?id=82149'
$ws1.Range('I45').NumberFormat = '@'
$ws1.Range('I45').Value = '//i0.hdslb.com/bfs/openplatform/202402/2odtsSXm1709101442352.jpeg'
$ws1.Range('F47').Value = 67
$ws1.Range('F48').Value = 538
$ws1.Range('F49').Value = 3035
$ws1.Range('F50').Value = 91
$ws2 = $wb.Worksheets.Item('演出')
$ws2.Range('F2').Value = 11
$ws2.Range('F4').Value = 19
$ws2.Range('F7').Value = 82
$ws2.Range('F15').Value = 14
$ws4 = $wb.Worksheets.Item('全部类型')
$ws4.Range('F3').Value = 5173
$ws4.Range('F4').Value = 5173
$ws4.Range('F5').Value = 122
$ws4.Range('F7').Value = 211
$ws4.Range('F9').Value = 213
$ws4.Range('F10').Value = 176
$ws4.Range('F11').Value = 8622
$ws4.Range('F12').Value = 8622
$ws4.Range('F13').Value = 275
$ws4.Range('F15').Value = 628
$ws4.Range('F16').Value = 2572
$ws4.Range('F17').Value = 19
$ws4.Range('F21').Value = 82
$ws4.Range('F26').Value = 6501
$ws4.Range('F27').Value = 203
$ws4.Range('F30').Value = 144
$ws4.Range('F33').Value = 7028
$ws4.Range('F37').Value = 109
$ws4.Range('F41').Value = 2542
$ws4.Range('F44').Value = 67
$ws4.Range('F45').Value = 538
$ws4.Range('F47').Value = 3036
$ws4.Range('F48').Value = 91
$ws4.Range('F49').Value = 14
